$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Append "请找dba索取" to the "输入网址：" line.
#    A trailing sentinel character is appended too so the insertion
#    point used for the relocated _GoBack bookmark below is never the
#    very last character position of the paragraph (that edge case
#    mis-places collapsed bookmarks in this runtime). The sentinel is
#    removed again right after the bookmark has been created.
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("输入网址：", $true, $false, $false, $false, $false, $true, 1, $false, "输入网址：请找dba索取~", 2)

# Locate the freshly written paragraph again to get reliable offsets.
$found = $d.Content
$found.Find.Execute("输入网址：请找dba索取~", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$sentinelStart = $found.End - 1

# ------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark (Word re-targets it automatically on
#    every edit) so it collapses right after the newly typed text,
#    matching where Word would leave it after typing there by hand.
#    Creating a bookmark with the same name removes the old one.
# ------------------------------------------------------------------
$bmRange = $d.Range($sentinelStart, $sentinelStart)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# 3. Drop the sentinel character now that the bookmark sits in the
#    right spot.
# ------------------------------------------------------------------
$sentinelRange = $d.Range($sentinelStart, $sentinelStart + 1)
$sentinelRange.Delete()
